$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.891.02"
$ws.Range("D3").Value = "2.667.39"
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'600.55"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'160.54"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("D7").Value = "'0.645"
$ws.Range("E7").Value = "  +4.51%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'5.89"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.401"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D13").Value = "'29.25"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "'0.0000197"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "3.147.71"
$ws.Range("D16").Value = "65.764.50"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "2.673.27"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "'12.63"
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").Value = "'356.49"
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'69.93"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("E24").Value = "  +9.83%  "
$ws.Range("D25").Value = "'0.0000114"
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("D26").Value = "'9.78"
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").Value = "'1.63"
$ws.Range("E27").Value = "  +2.72%  "
$ws.Range("D28").Value = "'571.43"
$ws.Range("E28").Value = "  +8.32%  "
$ws.Range("D29").Value = "'8.18"
$ws.Range("E29").Value = "  +2.05%  "
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "'2.15"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").Value = "'1.84"
$ws.Range("E33").Value = "  +3.98%  "
$ws.Range("E34").Value = "  +4.43%  "
$ws.Range("D35").Value = "'5.54"
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("D36").Value = "'0.425"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "'20.66"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").Value = "'1.99"
$ws.Range("E38").Value = "  +2.97%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "'154.29"
$ws.Range("E40").Value = "  -3.13%  "
$ws.Range("D41").Value = "'2.50"
$ws.Range("E41").Value = "  +7.17%  "
$ws.Range("D42").Value = "'162.74"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").Value = "'4.12"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").Value = "'0.0620"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("D45").Value = "'23.72"
$ws.Range("E45").Value = "  +4.33%  "
$ws.Range("D46").Value = "'0.646"
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("D49").Value = "'19.88"
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("E50").Value = "  -6.38%  "
$ws.Range("D51").Value = "'0.820"
$ws.Range("E51").Value = "  +1.26%  "
